$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the rows that changed per the diff.
$values = @{
    3  = -6
    4  = -1
    5  = 5
    6  = -3
    7  = -3
    8  = -1
    9  = 2
    10 = -5
    11 = -2
    12 = 1
    13 = -6
    14 = 2
    16 = -1
    18 = -2
    20 = 1
    21 = -1
    22 = -1
    23 = 3
    24 = -1
    25 = 3
    28 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
